# TOR-23 Fix sheet Band/Hero
# Strip the leading enum-type qualifiers ("ItemSlotTypeTOR.", "SkillTypeTOR.",
# "MagicItemType.") from the demo data on the "Items" sheet, row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Items")
$ws.Activate()

$ws.Range("E2").Value = "WEAPON"
$ws.Range("F2").Value = "BATTLE"
$ws.Range("G2").Value = "NONE"
$ws.Range("I2").Value = "UNUSUAL"

$ws.Range("I2").Select() | Out-Null
